# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14033
$wsExhibit.Range("F3").Value = 331
$wsExhibit.Range("F4").Value = 678
$wsExhibit.Range("F6").Value = 528
$wsExhibit.Range("F7").Value = 1456

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14033
$wsAll.Range("F3").Value = 331
$wsAll.Range("F4").Value = 678
$wsAll.Range("F8").Value = 528
$wsAll.Range("F9").Value = 1456
